$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1317.1578
$ws.Range("I15").Value = 1317.1578
$ws.Range("K15").Value = 3951.4734
$ws.Range("M15").Value = -3782.4734
$ws.Range("H19").Value = 3152.2727
$ws.Range("I19").Value = 5700
$ws.Range("J19").Value = 1388.4615
$ws.Range("K19").Value = 5700
$ws.Range("L19").Value = 1388.4615
$ws.Range("M19").Value = -5525
$ws.Range("N19").Value = -1738.4615
$ws.Range("H28").Value = 296
$ws.Range("I28").Value = 256
$ws.Range("K28").Value = 256
$ws.Range("M28").Value = 229
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""
$ws.Range("H64").Value = 2166.6667
$ws.Range("I64").Value = 2118.182
$ws.Range("J64").Value = 2300
$ws.Range("K64").Value = 2118.182
$ws.Range("L64").Value = 2300
$ws.Range("M64").Value = -1870.182
$ws.Range("N64").Value = -2796
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""
$ws.Range("H67").Value = 2166.6667
$ws.Range("I67").Value = 2118.182
$ws.Range("J67").Value = 2300
$ws.Range("K67").Value = 2118.182
$ws.Range("L67").Value = 2300
$ws.Range("M67").Value = -1260.182
$ws.Range("N67").Value = -4016
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 15000
$ws.Range("M69").Value = ""
$ws.Range("N69").Value = -16748
$ws.Range("H70").Value = 2254.3333
$ws.Range("I70").Value = 839.8
$ws.Range("J70").Value = 3264.7144
$ws.Range("K70").Value = 2519.4
$ws.Range("L70").Value = 9794.143199999999
$ws.Range("M70").Value = -2249.4
$ws.Range("N70").Value = -10334.1432
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 45000
$ws.Range("M72").Value = ""
$ws.Range("N72").Value = -53736
$ws.Range("H73").Value = 2254.3333
$ws.Range("I73").Value = 839.8
$ws.Range("J73").Value = 3264.7144
$ws.Range("K73").Value = 2519.4
$ws.Range("L73").Value = 9794.143199999999
$ws.Range("M73").Value = -1583.4
$ws.Range("N73").Value = -11666.1432
$ws.Range("H74").Value = 3705.3845
$ws.Range("I74").Value = 3680
$ws.Range("J74").Value = 3900
$ws.Range("K74").Value = 3680
$ws.Range("L74").Value = 3900
$ws.Range("M74").Value = -2744
$ws.Range("N74").Value = -5772
$ws.Range("H77").Value = 3705.3845
$ws.Range("I77").Value = 3680
$ws.Range("J77").Value = 3900
$ws.Range("K77").Value = 18400
$ws.Range("L77").Value = 19500
$ws.Range("M77").Value = -13720
$ws.Range("N77").Value = -28860
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -9508
$ws.Range("H116").Value = 2430
$ws.Range("I116").Value = 1980
$ws.Range("J116").Value = 2520
$ws.Range("K116").Value = 1980
$ws.Range("L116").Value = 2520
$ws.Range("M116").Value = 1462
$ws.Range("N116").Value = -9404
$ws.Range("H137").Value = 1648.5652
$ws.Range("I137").Value = 1495.4375
$ws.Range("J137").Value = 1998.5714
$ws.Range("K137").Value = 4486.3125
$ws.Range("L137").Value = 5995.7142
$ws.Range("M137").Value = -1936.3125
$ws.Range("N137").Value = -11095.7142

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1196.4706
$ws.Range("I2").Value = 1004.1
$ws.Range("J2").Value = 1471.2858
$ws.Range("K2").Value = 1004.1
$ws.Range("L2").Value = 1471.2858
$ws.Range("M2").Value = -891.1
$ws.Range("N2").Value = -1697.2858
$ws.Range("H45").Value = 1044.7693
$ws.Range("I45").Value = 774.6667
$ws.Range("J45").Value = 1652.5
$ws.Range("K45").Value = 774.6667
$ws.Range("L45").Value = 1652.5
$ws.Range("M45").Value = -397.6667
$ws.Range("N45").Value = -2406.5
$ws.Range("H61").Value = 1726.125
$ws.Range("I61").Value = 1119.4166
$ws.Range("J61").Value = 2332.8333
$ws.Range("K61").Value = 1119.4166
$ws.Range("L61").Value = 2332.8333
$ws.Range("M61").Value = -907.4166
$ws.Range("N61").Value = -2756.8333
$ws.Range("H110").Value = 1250.12
$ws.Range("I110").Value = 1111.35
$ws.Range("J110").Value = 1805.2
$ws.Range("K110").Value = 1111.35
$ws.Range("L110").Value = 1805.2
$ws.Range("M110").Value = 933.6500000000001
$ws.Range("N110").Value = -5895.2
$ws.Range("H116").Value = 1196.4706
$ws.Range("I116").Value = 1004.1
$ws.Range("J116").Value = 1471.2858
$ws.Range("K116").Value = 1004.1
$ws.Range("L116").Value = 1471.2858
$ws.Range("M116").Value = 1289.9
$ws.Range("N116").Value = -6059.2858
$ws.Range("H136").Value = 1726.125
$ws.Range("I136").Value = 1119.4166
$ws.Range("J136").Value = 2332.8333
$ws.Range("K136").Value = 3358.2498
$ws.Range("L136").Value = 6998.499899999999
$ws.Range("M136").Value = -808.2498000000001
$ws.Range("N136").Value = -12098.4999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1196.4706
$ws.Range("I3").Value = 1004.1
$ws.Range("J3").Value = 1471.2858
$ws.Range("K3").Value = 1004.1
$ws.Range("L3").Value = 1471.2858
$ws.Range("M3").Value = -890.1
$ws.Range("N3").Value = -1699.2858
$ws.Range("H63").Value = 35000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372
$ws.Range("H66").Value = 35000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864
$ws.Range("H107").Value = 1038.9231
$ws.Range("I107").Value = 568.6875
$ws.Range("J107").Value = 1791.3
$ws.Range("K107").Value = 568.6875
$ws.Range("L107").Value = 1791.3
$ws.Range("M107").Value = 1351.3125
$ws.Range("N107").Value = -5631.3

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1644.25
$ws.Range("I16").Value = 1110.9
$ws.Range("J16").Value = 2533.1667
$ws.Range("K16").Value = 1110.9
$ws.Range("L16").Value = 2533.1667
$ws.Range("M16").Value = -823.9000000000001
$ws.Range("N16").Value = -3107.1667
$ws.Range("H107").Value = 457.5
$ws.Range("I107").Value = 472.81818
$ws.Range("J107").Value = 373.25
$ws.Range("K107").Value = 472.81818
$ws.Range("L107").Value = 373.25
$ws.Range("M107").Value = 1447.18182
$ws.Range("N107").Value = -4213.25
$ws.Range("H113").Value = 1644.25
$ws.Range("I113").Value = 1110.9
$ws.Range("J113").Value = 2533.1667
$ws.Range("K113").Value = 1110.9
$ws.Range("L113").Value = 2533.1667
$ws.Range("M113").Value = 1059.1
$ws.Range("N113").Value = -6873.1667

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1405.95
$ws.Range("I102").Value = 1333.7778
$ws.Range("K102").Value = 1333.7778
$ws.Range("M102").Value = 288.2221999999999
$ws.Range("H107").Value = 517.0454999999999
$ws.Range("I107").Value = 448.8
$ws.Range("J107").Value = 663.2857
$ws.Range("K107").Value = 448.8
$ws.Range("L107").Value = 663.2857
$ws.Range("M107").Value = 1471.2
$ws.Range("N107").Value = -4503.2857
$ws.Range("H113").Value = 25003596
$ws.Range("I113").Value = 50005196
$ws.Range("J113").Value = 1996
$ws.Range("K113").Value = 50005196
$ws.Range("L113").Value = 1996
$ws.Range("M113").Value = -50003026
$ws.Range("N113").Value = -6336

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2367.9
$ws.Range("I61").Value = 1567.8
$ws.Range("J61").Value = 3168
$ws.Range("K61").Value = 1567.8
$ws.Range("L61").Value = 3168
$ws.Range("M61").Value = -1365.8
$ws.Range("N61").Value = -3572
$ws.Range("H113").Value = 2367.9
$ws.Range("I113").Value = 1567.8
$ws.Range("J113").Value = 3168
$ws.Range("K113").Value = 1567.8
$ws.Range("L113").Value = 3168
$ws.Range("M113").Value = 602.2
$ws.Range("N113").Value = -7508

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""
$ws.Range("H81").Value = 1549.9
$ws.Range("I81").Value = 1314.1428
$ws.Range("J81").Value = 2100
$ws.Range("K81").Value = 2628.2856
$ws.Range("L81").Value = 4200
$ws.Range("M81").Value = -1567.2856
$ws.Range("N81").Value = -6322
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""
$ws.Range("H84").Value = 1549.9
$ws.Range("I84").Value = 1314.1428
$ws.Range("J84").Value = 2100
$ws.Range("K84").Value = 13141.428
$ws.Range("L84").Value = 21000
$ws.Range("M84").Value = -7837.428
$ws.Range("N84").Value = -31608
$ws.Range("H107").Value = 1157.6471
$ws.Range("I107").Value = 1357.8462
$ws.Range("J107").Value = 507
$ws.Range("K107").Value = 4073.5386
$ws.Range("L107").Value = 1521
$ws.Range("M107").Value = -2153.5386
$ws.Range("N107").Value = -5361
$ws.Range("H113").Value = 503.05554
$ws.Range("I113").Value = 410.85715
$ws.Range("K113").Value = 1232.57145
$ws.Range("M113").Value = 937.4285500000001
